$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers for new columns I and J
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the header formatting (bold, border, centered) from H1 onto I1:J1
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Data values for I2:J19
$data = @(
    @(7, 7),
    @(6, 6),
    @(6, 6),
    @(5, 6),
    @(8, 8),
    @(6, 6),
    @(6, 6),
    @(8, 9),
    @(7, 7),
    @(8, 8),
    @(7, 7),
    @(6, 6),
    @(6, 6),
    @(7, 7),
    @(4, 5),
    @(5, 5),
    @(6, 6),
    @(5, 5)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 9).Value = $data[$i][0]
    $ws.Cells.Item($row, 10).Value = $data[$i][1]
}
